$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Status text: "In Translation" -> "Handed back: in sync with en-US"
#    (appears in the Status columns of Overview, zh-cn and de-de sheets)
# ------------------------------------------------------------------
foreach ($name in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($name)
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("In Translation" -eq $cell.Value2) {
            $cell.Value = "Handed back: in sync with en-US"
        }
    }
}

# ------------------------------------------------------------------
# 2. zh-cn sheet: fill in "Latest Target File" (I), "Latest Handback
#    File" (J) and "Latest Handback DateTime" (K) for both data rows,
#    plus a hyperlink on the new Latest Target File cell.
# ------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I2").Value = "53ca9609-2bc2-4999-b0f9-c2be6ce90c45.md"
$wsZh.Range("J2").Value = "53ca9609-2bc2-4999-b0f9-c2be6ce90c45.dc0388c8af5e1f2cf40a9b0af8ffb45c7a7fffda.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-23 20:25:26"

$wsZh.Range("I3").Value = "e07b38d4-f89b-408c-af98-b00b03f6fdd6.md"
$wsZh.Range("J3").Value = "e07b38d4-f89b-408c-af98-b00b03f6fdd6.06e1c968f384dc3c0b6c3c84a485203fa3078233.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-23 20:25:26"

# Re-create the hyperlinks in display order (A2, I2, A3, I3) so the
# relationship ids line up the same way Excel would renumber them.
$zhA2Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a910ed7f6532a5cb7d61e3b8199902d0c383326/e2e/53ca9609-2bc2-4999-b0f9-c2be6ce90c45.md"
$zhA3Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a910ed7f6532a5cb7d61e3b8199902d0c383326/e2e/e07b38d4-f89b-408c-af98-b00b03f6fdd6.md"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhA2Target, "", "", "53ca9609-2bc2-4999-b0f9-c2be6ce90c45.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhA2Target, "", "", "53ca9609-2bc2-4999-b0f9-c2be6ce90c45.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhA3Target, "", "", "e07b38d4-f89b-408c-af98-b00b03f6fdd6.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhA3Target, "", "", "e07b38d4-f89b-408c-af98-b00b03f6fdd6.md")

# ------------------------------------------------------------------
# 3. de-de sheet: same treatment.
# ------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I2").Value = "53ca9609-2bc2-4999-b0f9-c2be6ce90c45.md"
$wsDe.Range("J2").Value = "53ca9609-2bc2-4999-b0f9-c2be6ce90c45.dc0388c8af5e1f2cf40a9b0af8ffb45c7a7fffda.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-23 20:25:33"

$wsDe.Range("I3").Value = "e07b38d4-f89b-408c-af98-b00b03f6fdd6.md"
$wsDe.Range("J3").Value = "e07b38d4-f89b-408c-af98-b00b03f6fdd6.06e1c968f384dc3c0b6c3c84a485203fa3078233.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-23 20:25:33"

$deA2Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a910ed7f6532a5cb7d61e3b8199902d0c383326/e2e/53ca9609-2bc2-4999-b0f9-c2be6ce90c45.md"
$deA3Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a910ed7f6532a5cb7d61e3b8199902d0c383326/e2e/e07b38d4-f89b-408c-af98-b00b03f6fdd6.md"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deA2Target, "", "", "53ca9609-2bc2-4999-b0f9-c2be6ce90c45.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deA2Target, "", "", "53ca9609-2bc2-4999-b0f9-c2be6ce90c45.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deA3Target, "", "", "e07b38d4-f89b-408c-af98-b00b03f6fdd6.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deA3Target, "", "", "e07b38d4-f89b-408c-af98-b00b03f6fdd6.md")

# ------------------------------------------------------------------
# 4. Column widths widened to fit the newly-populated columns.
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 29.9777047293527

$wsZh.Range("C1").EntireColumn.ColumnWidth = 29.9777047293527
$wsZh.Range("I1:J1").EntireColumn.ColumnWidth = 40

$wsDe.Range("C1").EntireColumn.ColumnWidth = 29.9777047293527
$wsDe.Range("I1:J1").EntireColumn.ColumnWidth = 40

Write-Host "Handback report generated."
